$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: replace the old time slot "12:25-12:30" with a new one "18:55-19:0"
# (the "Посыл / Заповедь" value in B6 is left as-is, same recurring text)
$ws.Range("C6").Value = "18:55-19:0"

# Row 7: replace the old time slot "12:50-12:55" with "19:0-19:5", and switch
# the B7 "Посыл / Заповедь" to the alternate (shorter) recurring text, matching
# the alternating pattern already used in rows 3 and 5
$ws.Range("B7").Value = $ws.Range("B5").Value2
$ws.Range("C7").Value = "19:0-19:5"

# Remove the trailing rows 8-11 (old 14:10-14:15 / 14:15-14:20 / 22:45-22:50 /
# 22:50-22:55 slots) entirely - the sheet now ends at row 7
$ws.Rows("8:11").Delete()

# Move the active selection, matching where the author left the cursor
$ws.Range("B13").Select()
